# Update NATMI TPM figures for the Mdk -> Itgb1 ligand/receptor sheet.
#
# Columns G,H,I,J (ligand average/total expression + derived specificity)
# depend only on the "Sending cluster" (column A); columns M,N,O,P
# (receptor average/total expression + derived specificity) depend only
# on the "Target cluster" (column D); columns Q,R,S,T (edge weights /
# derived specificities) are recomputed from the new ligand & receptor
# TPM figures. Rows 2-17 hold one row per Sending x Target cluster pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row number -> new values for columns G,H,I,J,M,N,O,P,Q,R,S,T (in that order).
$newValues = @{
  2  = @(2.180165333333334, 6.540496, 0.01970539991828544, 0.01970539991828544, 145.7007446666667, 437.1022340000001, 0.2865937750105843, 0.2865937750105843, 317.6517125631183, 2858.865413068064, 0.005647444950674681, 0.005647444950674681)
  3  = @(2.180165333333334, 6.540496, 0.01970539991828544, 0.01970539991828544, 168.7997026666667, 506.3991080000001, 0.3320294904365841, 0.3320294904365841, 368.011260030841, 3312.101340277568, 0.00654277389371742, 0.00654277389371742)
  4  = @(2.180165333333334, 6.540496, 0.01970539991828544, 0.01970539991828544, 128.1261546666667, 384.378464, 0.2520245069956105, 0.2520245069956105, 279.3362006975716, 2514.025806278144, 0.004966243699557229, 0.004966243699557229)
  5  = @(2.180165333333334, 6.540496, 0.01970539991828544, 0.01970539991828544, 65.761079, 197.283237, 0.1293522275572212, 0.1293522275572212, 143.3700247183947, 1290.330222465552, 0.002548937374336106, 0.002548937374336106)
  6  = @(81.17653533333333, 243.529606, 0.733713204346044, 0.7337132043460441, 145.7007446666667, 437.1022340000001, 0.2865937750105843, 0.2865937750105843, 11827.48164752665, 106447.3348277398, 0.210277637008645, 0.210277637008645)
  7  = @(81.17653533333333, 243.529606, 0.733713204346044, 0.7337132043460441, 168.7997026666667, 506.3991080000001, 0.3320294904365841, 0.3320294904365841, 13702.57502777683, 123323.1752499915, 0.2436144213656103, 0.2436144213656104)
  8  = @(81.17653533333333, 243.529606, 0.733713204346044, 0.7337132043460441, 128.1261546666667, 384.378464, 0.2520245069956105, 0.2520245069956105, 10400.8373214228, 93607.53589280519, 0.1849137086014813, 0.1849137086014814)
  9  = @(81.17653533333333, 243.529606, 0.733713204346044, 0.7337132043460441, 65.761079, 197.283237, 0.1293522275572212, 0.1293522275572212, 5338.256553001625, 48044.30897701462, 0.09490743737030745, 0.09490743737030746)
  10 = @(25.672264, 77.016792, 0.2320384702908474, 0.2320384702908474, 145.7007446666667, 437.1022340000001, 0.2865937750105843, 0.2865937750105843, 3740.467982079259, 33664.21183871333, 0.06650078114833526, 0.06650078114833526)
  11 = @(25.672264, 77.016792, 0.2320384702908474, 0.2320384702908474, 168.7997026666667, 506.3991080000001, 0.3320294904365841, 0.3320294904365841, 4333.470529980171, 39001.23476982154, 0.07704361505235453, 0.07704361505235453)
  12 = @(25.672264, 77.016792, 0.2320384702908474, 0.2320384702908474, 128.1261546666667, 384.378464, 0.2520245069956105, 0.2520245069956105, 3289.288467907499, 29603.59621116749, 0.05847938107906642, 0.05847938107906642)
  13 = @(25.672264, 77.016792, 0.2320384702908474, 0.2320384702908474, 65.761079, 197.283237, 0.1293522275572212, 0.1293522275572212, 1688.235781012856, 15194.1220291157, 0.03001469301109121, 0.03001469301109121)
  14 = @(1.608999666666667, 4.826999, 0.01454292544482312, 0.01454292544482312, 145.7007446666667, 437.1022340000001, 0.2865937750105843, 0.2865937750105843, 234.4324496017518, 2109.892046415766, 0.00416791190292934, 0.00416791190292934)
  15 = @(1.608999666666667, 4.826999, 0.01454292544482312, 0.01454292544482312, 168.7997026666667, 506.3991080000001, 0.3320294904365841, 0.3320294904365841, 271.5986653240992, 2444.387987916892, 0.004828680124901855, 0.004828680124901855)
  16 = @(1.608999666666667, 4.826999, 0.01454292544482312, 0.01454292544482312, 128.1261546666667, 384.378464, 0.2520245069956105, 0.2520245069956105, 206.1549401499485, 1855.394461349536, 0.003665173615505467, 0.003665173615505467)
  17 = @(1.608999666666667, 4.826999, 0.01454292544482312, 0.01454292544482312, 65.761079, 197.283237, 0.1293522275572212, 0.1293522275572212, 105.8095541906403, 952.2859877157629, 0.001881159801486463, 0.001881159801486463)
}

$cols = @(7, 8, 9, 10, 13, 14, 15, 16, 17, 18, 19, 20)

for ($r = 2; $r -le 17; $r++) {
  $vals = $newValues[$r]
  for ($i = 0; $i -lt $cols.Count; $i++) {
    $ws.Cells.Item($r, $cols[$i]).Value = $vals[$i]
  }
}
